$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data (hour label, count) for rows 2..25, in the new order/values. This is
# also the exact order in which the labels must be (re)written so the
# rebuilt shared-string table ends up in the same insertion order as the
# target workbook.
$data = @(
    @("09", 81),
    @("08", 71),
    @("16", 59),
    @("11", 58),
    @("13", 58),
    @("10", 55),
    @("24", 54),
    @("14", 47),
    @("17", 42),
    @("07", 40),
    @("15", 38),
    @("18", 28),
    @("19", 28),
    @("20", 26),
    @("21", 26),
    @("06", 19),
    @("05", 17),
    @("22", 16),
    @("12", 14),
    @("04", 13),
    @("23", 9),
    @("01", 9),
    @("02", 7),
    @("03", 6)
)

# Clear out the existing label/count rows first (including the row that
# will become new row 25). This drops the old shared-string entries so a
# fresh shared-string table gets rebuilt -- in the exact insertion order
# used below -- once the new labels are written.
$ws.Range("A2:B25").Value = $null

# Format the label column as text up front so values such as "09"/"08"
# are stored as text (preserving the leading zero) instead of being
# auto-converted to the numbers 9/8.
$ws.Range("A2:A25").NumberFormat = "@"

$row = 2
foreach ($pair in $data) {
    $label = $pair[0]
    $count = $pair[1]

    $ws.Cells.Item($row, 1).Value = $label
    $ws.Cells.Item($row, 2).Value = $count

    $row = $row + 1
}
